$wb = $excel.ActiveWorkbook

$newTimestamp = "February 03 2026 18.05.36 EST"

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A2").Value = "Version: Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newTimestamp)"
$wsAbout.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Wangjiayu Coal Mine, China, M2184, version 'Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newTimestamp)'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

$wsData = $wb.Worksheets.Item("Boundaries and methane sources")
$newVersionString = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newTimestamp)"

for ($row = 2; $row -le 7; $row++) {
    $wsData.Range("S$row").Value = $newVersionString
}
